# Updates the cryptos worksheet with refreshed prices / 1h volume percentages,
# and replaces the last row's coin (FirstDigitalUSD -> Monero).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) updates. These are written as text (matching the sheet's
# existing inline-string "Price" column) even though many of them look like
# plain numbers, so we force the cell to Text format first, then restore the
# default "Normal" style so we don't leave a different number format applied
# to the cell versus the rest of the sheet.
$priceUpdates = @(
    @{Cell="D2"; Value="70.678.75"},
    @{Cell="D3"; Value="3.628.37"},
    @{Cell="D5"; Value="594.21"},
    @{Cell="D6"; Value="192.28"},
    @{Cell="D7"; Value="0.658"},
    @{Cell="D8"; Value="3.611.89"},
    @{Cell="D11"; Value="0.666"},
    @{Cell="D12"; Value="58.12"},
    @{Cell="D13"; Value="0.0000296"},
    @{Cell="D14"; Value="9.79"},
    @{Cell="D15"; Value="4.209.12"},
    @{Cell="D17"; Value="3.626.96"},
    @{Cell="D18"; Value="70.615.02"},
    @{Cell="D22"; Value="495.62"},
    @{Cell="D24"; Value="17.08"},
    @{Cell="D25"; Value="4.49"},
    @{Cell="D26"; Value="91.26"},
    @{Cell="D27"; Value="3.13"},
    @{Cell="D28"; Value="11.29"},
    @{Cell="D29"; Value="9.48"},
    @{Cell="D30"; Value="32.39"},
    @{Cell="D31"; Value="7.61"},
    @{Cell="D32"; Value="12.29"},
    @{Cell="D33"; Value="618.94"},
    @{Cell="D35"; Value="65.29"},
    @{Cell="D39"; Value="38.16"},
    @{Cell="D41"; Value="3.64"},
    @{Cell="D42"; Value="3.343.81"},
    @{Cell="D44"; Value="0.0450"},
    @{Cell="D46"; Value="3.40"},
    @{Cell="D48"; Value="9.23"},
    @{Cell="D50"; Value="3.37"}
)

foreach ($item in $priceUpdates) {
    $cell = $ws.Range($item.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $item.Value
    $cell.Style = "Normal"
}

# Volume(1h) (column E) updates - these already contain non-numeric
# characters (%, spaces) so Excel keeps them as text automatically.
$volumeUpdates = @(
    @{Cell="E2"; Value="  +7.28%  "},
    @{Cell="E3"; Value="  +7.11%  "},
    @{Cell="E4"; Value="  +0.07%  "},
    @{Cell="E5"; Value="  +5.22%  "},
    @{Cell="E6"; Value="  +8.85%  "},
    @{Cell="E7"; Value="  +4.39%  "},
    @{Cell="E8"; Value="  +6.87%  "},
    @{Cell="E9"; Value="  -0.06%  "},
    @{Cell="E10"; Value="  +3.05%  "},
    @{Cell="E11"; Value="  +5.04%  "},
    @{Cell="E12"; Value="  +7.72%  "},
    @{Cell="E13"; Value="  +6.05%  "},
    @{Cell="E14"; Value="  +5.69%  "},
    @{Cell="E16"; Value="  +6.84%  "},
    @{Cell="E17"; Value="  +7.43%  "},
    @{Cell="E18"; Value="  +7.46%  "},
    @{Cell="E19"; Value="  +6.25%  "},
    @{Cell="E20"; Value="  +0.95%  "},
    @{Cell="E21"; Value="  +5.65%  "},
    @{Cell="E22"; Value="  +6.74%  "},
    @{Cell="E23"; Value="  +10.97%  "},
    @{Cell="E24"; Value="  +15.50%  "},
    @{Cell="E25"; Value="  +9.27%  "},
    @{Cell="E26"; Value="  +1.98%  "},
    @{Cell="E27"; Value="  +6.91%  "},
    @{Cell="E28"; Value="  +5.84%  "},
    @{Cell="E29"; Value="  +8.59%  "},
    @{Cell="E30"; Value="  +4.06%  "},
    @{Cell="E31"; Value="  +15.12%  "},
    @{Cell="E32"; Value="  +6.90%  "},
    @{Cell="E33"; Value="  +6.55%  "},
    @{Cell="E34"; Value="  +8.95%  "},
    @{Cell="E35"; Value="  +4.76%  "},
    @{Cell="E36"; Value="  +10.96%  "},
    @{Cell="E37"; Value="  +7.73%  "},
    @{Cell="E38"; Value="  +3.91%  "},
    @{Cell="E39"; Value="  +5.78%  "},
    @{Cell="E40"; Value="  +0.00%  "},
    @{Cell="E41"; Value="  +1.27%  "},
    @{Cell="E42"; Value="  +7.72%  "},
    @{Cell="E43"; Value="  +8.20%  "},
    @{Cell="E44"; Value="  +7.57%  "},
    @{Cell="E45"; Value="  +10.30%  "},
    @{Cell="E46"; Value="  +7.28%  "},
    @{Cell="E47"; Value="  +3.70%  "},
    @{Cell="E48"; Value="  +8.71%  "},
    @{Cell="E49"; Value="  +6.80%  "},
    @{Cell="E50"; Value="  +6.07%  "}
)

foreach ($item in $volumeUpdates) {
    $ws.Range($item.Cell).Value = $item.Value
}

# Row 51: FirstDigitalUSD -> Monero
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"

$d51 = $ws.Range("D51")
$d51.NumberFormat = "@"
$d51.Value = "143.37"
$d51.Style = "Normal"

$ws.Range("E51").Value = "  +1.81%  "
